# Handback status report: record the actual handback completion
# timestamps for row 3 (the f86032e9... file) on both locale sheets.
# Previously these cells were placeholders copied from row 2; now each
# locale/file pair gets its own "Correspond Handback DateTime" /
# completion datetime values.

$wb = $excel.ActiveWorkbook

# zh-cn sheet (row 3: f86032e9-...zh-cn.xlf)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-12 00:32:53"
$wsZhCn.Range("H3").Value = "2016-03-12 00:33:10"

# de-de sheet (row 3: f86032e9-...de-de.xlf)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-12 00:32:56"
$wsDeDe.Range("H3").Value = "2016-03-12 00:33:15"
